# GDP Growth Rates.xlsx update — refresh STEO source data from "July STEO"
# to "September STEO" (per commit message: most recent export), update the
# Real GDP figures for 2020/2021, and refresh the related narrative text
# on the About sheet. Also turn off iterative calculation and nudge the
# saved cell-selections to match the author's last editing position.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Data")

# --- Data sheet: refresh STEO label + the raw Real GDP figures -------------
$wsData.Range("A3").Value = "September STEO"
$wsData.Range("C3").Value = 18168
$wsData.Range("D3").Value = 18726

# --- About sheet: refresh the source-citation text -------------------------
$wsAbout.Range("B6").Value  = "January 2020 and September 2020"
$wsAbout.Range("A28").Value = "SARS-CoV-2 pandemic.  It uses the latest data available as of September 9,"

# --- Turn off iterative calculation (workbook no longer needs it) ----------
$excel.Iteration = $false

# --- Restore the cell selections left by the editor on save ----------------
# (About stays the active/visible tab, so select it last.)
$wsData.Range("D4").Select() | Out-Null
$wsAbout.Range("A29").Select() | Out-Null
